$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 883.3333
$ws.Range("I28").Value = 860
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 860
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -375
$ws.Range("N28").Value = -1970
$ws.Range("H62").Value = 3329.3333
$ws.Range("I62").Value = 3067.1428
$ws.Range("K62").Value = 3067.1428
$ws.Range("M62").Value = -2443.1428
$ws.Range("H65").Value = 3329.3333
$ws.Range("I65").Value = 3067.1428
$ws.Range("K65").Value = 15335.714
$ws.Range("M65").Value = -12215.714
$ws.Range("H113").Value = 196323.12
$ws.Range("I113").Value = 11996.777
$ws.Range("J113").Value = 433314.16
$ws.Range("K113").Value = 11996.777
$ws.Range("L113").Value = 433314.16
$ws.Range("M113").Value = -8742.777
$ws.Range("N113").Value = -439822.16
$ws.Range("H121").Value = 3180.4
$ws.Range("J121").Value = 3180.4
$ws.Range("L121").Value = 9541.200000000001
$ws.Range("N121").Value = -13035.2
$ws.Range("H141").Value = 3541.56
$ws.Range("I141").Value = 3501.6956
$ws.Range("K141").Value = 10505.0868
$ws.Range("M141").Value = -5325.086800000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2269.7083
$ws.Range("I2").Value = 2654.4736
$ws.Range("K2").Value = 2654.4736
$ws.Range("M2").Value = -2541.4736
$ws.Range("H31").Value = 6583
$ws.Range("I31").Value = 6583
$ws.Range("K31").Value = 6583
$ws.Range("M31").Value = -6289
$ws.Range("H63").Value = 3908.3333
$ws.Range("I63").Value = 3908.3333
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3908.3333
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3222.3333
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 3908.3333
$ws.Range("I66").Value = 3908.3333
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 19541.6665
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -16109.6665
$ws.Range("N66").ClearContents()
$ws.Range("H116").Value = 2269.7083
$ws.Range("I116").Value = 2654.4736
$ws.Range("K116").Value = 2654.4736
$ws.Range("M116").Value = -360.4735999999998
$ws.Range("H122").Value = 1085.9678
$ws.Range("J122").Value = 1243.2
$ws.Range("L122").Value = 3729.6
$ws.Range("N122").Value = -8629.6
$ws.Range("H132").Value = 3487430.2
$ws.Range("I132").Value = 1755.3448
$ws.Range("K132").Value = 5266.0344
$ws.Range("M132").Value = -2736.0344

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2269.7083
$ws.Range("I3").Value = 2654.4736
$ws.Range("K3").Value = 2654.4736
$ws.Range("M3").Value = -2540.4736
$ws.Range("H86").Value = 3184.7097
$ws.Range("I86").Value = 3286.3845
$ws.Range("K86").Value = 3286.3845
$ws.Range("M86").Value = -2163.3845
$ws.Range("H89").Value = 3184.7097
$ws.Range("I89").Value = 3286.3845
$ws.Range("K89").Value = 16431.9225
$ws.Range("M89").Value = -10815.9225
$ws.Range("H102").Value = 11396.091
$ws.Range("I102").Value = 11396.091
$ws.Range("K102").Value = 11396.091
$ws.Range("M102").Value = -8151.091
$ws.Range("H105").Value = 1495.2667
$ws.Range("I105").Value = 1409.625
$ws.Range("K105").Value = 1409.625
$ws.Range("M105").Value = 337.375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 101589.23
$ws.Range("I31").Value = 201154.1
$ws.Range("J31").Value = 18618.5
$ws.Range("K31").Value = 201154.1
$ws.Range("L31").Value = 18618.5
$ws.Range("M31").Value = -200859.1
$ws.Range("N31").Value = -19208.5
$ws.Range("H34").Value = 101589.23
$ws.Range("I34").Value = 201154.1
$ws.Range("J34").Value = 18618.5
$ws.Range("K34").Value = 201154.1
$ws.Range("L34").Value = 18618.5
$ws.Range("M34").Value = -200952.1
$ws.Range("N34").Value = -19022.5
$ws.Range("H58").Value = 15403.357
$ws.Range("I58").Value = 5250.08
$ws.Range("K58").Value = 5250.08
$ws.Range("M58").Value = -5047.08
$ws.Range("H134").Value = 31255288
$ws.Range("I134").Value = 1970.0416
$ws.Range("K134").Value = 5910.1248
$ws.Range("M134").Value = -3375.1248
$ws.Range("H136").Value = 15403.357
$ws.Range("I136").Value = 5250.08
$ws.Range("K136").Value = 15750.24
$ws.Range("M136").Value = -13200.24

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 4999.6
$ws.Range("I31").Value = 3749.75
$ws.Range("J31").Value = 9999
$ws.Range("K31").Value = 11249.25
$ws.Range("L31").Value = 29997
$ws.Range("M31").Value = -10961.25
$ws.Range("N31").Value = -30573
$ws.Range("H107").Value = 963.75
$ws.Range("J107").Value = 1131.25
$ws.Range("L107").Value = 3393.75
$ws.Range("N107").Value = -7233.75
$ws.Range("H131").Value = 1477.12
$ws.Range("I131").Value = 1040
$ws.Range("J131").Value = 1500.1263
$ws.Range("K131").Value = 3120
$ws.Range("L131").Value = 4500.3789
$ws.Range("M131").Value = 1920
$ws.Range("N131").Value = -14580.3789
$ws.Range("H132").Value = 3368273.8
$ws.Range("I132").Value = 1348
$ws.Range("J132").Value = 5051736.5
$ws.Range("K132").Value = 12132
$ws.Range("L132").Value = 45465628.5
$ws.Range("M132").Value = -9602
$ws.Range("N132").Value = -45470688.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 70.5
$ws.Range("I2").Value = 79.82353000000001
$ws.Range("K2").Value = 79.82353000000001
$ws.Range("M2").Value = 33.17646999999999
$ws.Range("H62").Value = 19000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 19000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 18095
$ws.Range("I70").Value = 15000
$ws.Range("J70").Value = 20416.25
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 20416.25
$ws.Range("M70").Value = -14730
$ws.Range("N70").Value = -20956.25
$ws.Range("H73").Value = 18095
$ws.Range("I73").Value = 15000
$ws.Range("J73").Value = 20416.25
$ws.Range("K73").Value = 15000
$ws.Range("L73").Value = 20416.25
$ws.Range("M73").Value = -14064
$ws.Range("N73").Value = -22288.25
$ws.Range("H97").Value = 1007.7879
$ws.Range("I97").Value = 953.9048
$ws.Range("J97").Value = 1102.0834
$ws.Range("K97").Value = 953.9048
$ws.Range("L97").Value = 1102.0834
$ws.Range("M97").Value = -457.9048
$ws.Range("N97").Value = -2094.0834
$ws.Range("H102").Value = 9355.857
$ws.Range("I102").Value = 10841.182
$ws.Range("J102").Value = 3909.6667
$ws.Range("K102").Value = 10841.182
$ws.Range("L102").Value = 3909.6667
$ws.Range("M102").Value = -9219.182000000001
$ws.Range("N102").Value = -7153.6667
$ws.Range("H113").Value = 3708.9443
$ws.Range("I113").Value = 2064.7778
$ws.Range("K113").Value = 2064.7778
$ws.Range("M113").Value = 105.2222000000002
$ws.Range("H122").Value = 2304.7307
$ws.Range("I122").Value = 2134.65
$ws.Range("K122").Value = 6403.950000000001
$ws.Range("M122").Value = -3953.950000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 34666.332
$ws.Range("J14").Value = 100000
$ws.Range("L14").Value = 100000
$ws.Range("N14").Value = -100344
$ws.Range("H47").Value = 33567.25
$ws.Range("I47").Value = 40059
$ws.Range("J47").Value = 31403.334
$ws.Range("K47").Value = 40059
$ws.Range("L47").Value = 31403.334
$ws.Range("M47").Value = -39569
$ws.Range("N47").Value = -32383.334
$ws.Range("H52").Value = 33567.25
$ws.Range("I52").Value = 40059
$ws.Range("J52").Value = 31403.334
$ws.Range("K52").Value = 40059
$ws.Range("L52").Value = 31403.334
$ws.Range("M52").Value = -39826
$ws.Range("N52").Value = -31869.334
$ws.Range("H61").Value = 5117.6875
$ws.Range("I61").Value = 5417.273
$ws.Range("J61").Value = 4458.6
$ws.Range("K61").Value = 5417.273
$ws.Range("L61").Value = 4458.6
$ws.Range("M61").Value = -5215.273
$ws.Range("N61").Value = -4862.6
$ws.Range("H113").Value = 5117.6875
$ws.Range("I113").Value = 5417.273
$ws.Range("J113").Value = 4458.6
$ws.Range("K113").Value = 5417.273
$ws.Range("L113").Value = 4458.6
$ws.Range("M113").Value = -3247.273
$ws.Range("N113").Value = -8798.6
$ws.Range("H122").Value = 5769.2915
$ws.Range("I122").Value = 5641.095
$ws.Range("K122").Value = 16923.285
$ws.Range("M122").Value = -14473.285

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 234687.5
$ws.Range("I9").Value = 234687.5
$ws.Range("K9").Value = 234687.5
$ws.Range("M9").Value = -234547.5
$ws.Range("H14").Value = 26252
$ws.Range("J14").Value = 51000
$ws.Range("L14").Value = 51000
$ws.Range("N14").Value = -51336
$ws.Range("H81").Value = 5993.3335
$ws.Range("I81").Value = 5993.3335
$ws.Range("K81").Value = 11986.667
$ws.Range("M81").Value = -10925.667
$ws.Range("H84").Value = 5993.3335
$ws.Range("I84").Value = 5993.3335
$ws.Range("K84").Value = 59933.335
$ws.Range("M84").Value = -54629.335
$ws.Range("H107").Value = 1140.375
$ws.Range("I107").Value = 1149.7333
$ws.Range("K107").Value = 3449.199900000001
$ws.Range("M107").Value = -1529.199900000001
$ws.Range("H119").Value = 223333.33
$ws.Range("J119").Value = 223333.33
$ws.Range("L119").Value = 223333.33
$ws.Range("N119").Value = -233009.33
$ws.Range("H122").Value = 3142.9375
$ws.Range("I122").Value = 2685.8
$ws.Range("K122").Value = 8057.400000000001
$ws.Range("M122").Value = -5607.400000000001
